$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "Trade off performance with risk") {
        $shp.Delete()
    }
}
